$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

# Insert a new row at row 6, shifting existing rows 6-14 down to 7-15
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the DADA2Tax data
$ws.Cells.Item(6, 1).Value = "DADA2Tax"
$ws.Cells.Item(6, 2).Value = "100 Australian species"
$ws.Cells.Item(6, 3).Value = "COI"
$ws.Cells.Item(6, 4).Value = 0.2564102564102564
$ws.Cells.Item(6, 5).Value = 0.1818181818181818
$ws.Cells.Item(6, 6).Value = 0.2127659574468085
$ws.Cells.Item(6, 7).Value = 0.2369668246445497
$ws.Cells.Item(6, 8).Value = 0.2525252525252525
